# "added 4wk low sales check" - update forecast comparison figures and summary totals

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Row 2 (W10)
$ws1.Range("H2").Value = 33.47
$ws1.Range("L2").Value = 0.98

# Row 3 (W11)
$ws1.Range("D3").Value = 7
$ws1.Range("H3").Value = 33.4
$ws1.Range("L3").Value = 1.07

# Row 4 (W12)
$ws1.Range("H4").Value = 33.35
$ws1.Range("L4").Value = 0.93

# Row 5 (W13)
$ws1.Range("H5").Value = 31.43
$ws1.Range("L5").Value = 1.2

# Row 6 (W14)
$ws1.Range("D6").Value = 6
$ws1.Range("H6").Value = 31.32
$ws1.Range("L6").Value = 0.86

# Row 7 (W15)
$ws1.Range("D7").Value = 7
$ws1.Range("H7").Value = 29.46
$ws1.Range("L7").Value = 1.2

# Row 8 (W16)
$ws1.Range("D8").Value = 6
$ws1.Range("H8").Value = 29.29
$ws1.Range("L8").Value = 0.96

# Row 9 (W17)
$ws1.Range("D9").Value = 6
$ws1.Range("H9").Value = 28.29
$ws1.Range("L9").Value = 1.01

# Row 10 (W18)
$ws1.Range("H10").Value = 27.29
$ws1.Range("L10").Value = 0.98

# Row 11 (W19)
$ws1.Range("D11").Value = 6
$ws1.Range("H11").Value = 26.29
$ws1.Range("L11").Value = 1.01

# Row 12 (W20)
$ws1.Range("H12").Value = 24.57
$ws1.Range("L12").Value = 1.09

# Row 13 (W21)
$ws1.Range("D13").Value = 6
$ws1.Range("H13").Value = 24.26
$ws1.Range("L13").Value = 0.87

# Row 14 (W22)
$ws1.Range("H14").Value = 22.6
$ws1.Range("L14").Value = 1.02

# Row 15 (W23)
$ws1.Range("D15").Value = 7
$ws1.Range("H15").Value = 21.6
$ws1.Range("L15").Value = 0.91

# Row 16 (W24)
$ws1.Range("D16").Value = 7
$ws1.Range("H16").Value = 20.6
$ws1.Range("L16").Value = 0.85

# Row 17 (W25)
$ws1.Range("H17").Value = 20.18
$ws1.Range("L17").Value = 0.98

# Summary sheet totals
$ws2.Range("B9").Value = "111"
$ws2.Range("B10").Value = "55"
$ws2.Range("B11").Value = "28"
